# Update column F (dSF) values for a set of rows based on repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 2
    10 = 1
    12 = 0
    22 = 4
    26 = 4
    30 = 1
    32 = 3
    40 = -1
    41 = 0
    44 = 2
    54 = -4
    55 = 0
    60 = -1
    68 = -1
    69 = -1
    70 = 5
    71 = -1
    73 = -4
    77 = -3
    81 = 7
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
